# Refresh the cryptocurrency price/volume snapshot (GitHub Actions bot update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "59.346.56"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.73%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.605.78"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("E4").Value = "  +0.27%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "539.82"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.33%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "141.00"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("E7").Value = "  +0.52%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.567"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.16%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "6.45"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.65%  "

$ws.Range("E10").Value = "  +1.41%  "

$ws.Range("E11").Value = "  +0.56%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.136"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.93%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "3.066.27"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "59.263.62"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "20.57"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.60%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.629.48"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("E17").Value = "  +0.43%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "343.23"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.14%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.34"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.27%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "10.11"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.90%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.42"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.85%  "

$ws.Range("E22").Value = "  +0.32%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "67.37"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.66%  "

$ws.Range("E24").Value = "  -0.95%  "

$ws.Range("E25").Value = "  +0.63%  "

$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("E27").Value = "  +1.52%  "

$ws.Range("E28").Value = "  +0.21%  "

$ws.Range("E29").Value = "  +1.34%  "

$ws.Range("E30").Value = "  +5.77%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "5.83"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.24%  "

$ws.Range("E32").Value = "  -0.28%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "149.44"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.31%  "

$ws.Range("E34").Value = "  -1.12%  "

$ws.Range("E35").Value = "  -1.74%  "

$ws.Range("E36").Value = "  +1.65%  "

$ws.Range("E37").Value = "  -0.47%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.831"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.24%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.812"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.53%  "

$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("E41").Value = "  +0.12%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "273.63"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.29%  "

$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "10.75"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.596"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("E46").Value = "  +0.26%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.944.98"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.09%  "

$ws.Range("E48").Value = "  +0.67%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "4.52"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.28%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "18.26"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.70%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "110.92"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.49%  "

